$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("Week 5")
$ws6 = $wb.Worksheets.Item("Week 6")

# D5: number of work log entries for week 6
$ws6.Range("D5").Value = 6

# --- Fill in time-log entries for Week 6 ---
$ws6.Range("A15").Value = 'Project management (update logs, weekly summary, meeting minutes, etc)'
$ws6.Range("B15").Value = 'Project Mangaement'
$ws6.Range("C15").Value = 2.75
$ws6.Range("A16").Value = 'Stand-up / project planning discussion'
$ws6.Range("B16").Value = 'Internal Meeting'
$ws6.Range("C16").Value = 1.25
$ws6.Range("A17").Value = 'Work on report'
$ws6.Range("B17").Value = 'Report'
$ws6.Range("C17").Value = 4.75
$ws6.Range("A26").Value = 'Work on report'
$ws6.Range("B26").Value = 'Report'
$ws6.Range("C26").Value = 6.25
$ws6.Range("A27").Value = 'Stand-up and go over status presentation'
$ws6.Range("B27").Value = 'Internal Meeting'
$ws6.Range("C27").Value = 0.5
$ws6.Range("A28").Value = 'Weekly status meeting'
$ws6.Range("B28").Value = 'UBC Meeting'
$ws6.Range("C28").Value = 0.5
$ws6.Range("A29").Value = 'Work on next week status presentation'
$ws6.Range("B29").Value = 'Project Documents'
$ws6.Range("C29").Value = 0.5
$ws6.Range("A30").Value = 'Setup this week''s sprint planning doc'
$ws6.Range("B30").Value = 'Project Documents'
$ws6.Range("C30").Value = 0.5
$ws6.Range("A37").Value = 'Work on report'
$ws6.Range("B37").Value = 'Report'
$ws6.Range("C37").Value = 3.25
$ws6.Range("A38").Value = 'Discuss report with Ryan'
$ws6.Range("B38").Value = 'Report'
$ws6.Range("C38").Value = 0.5
$ws6.Range("A39").Value = 'Discuss dashboard with Mitch'
$ws6.Range("B39").Value = 'Dashboard'
$ws6.Range("C39").Value = 1
$ws6.Range("A40").Value = 'Review test environment'
$ws6.Range("B40").Value = 'Test Environment Example'
$ws6.Range("C40").Value = 3
$ws6.Range("A41").Value = 'Update sprint planning doc'
$ws6.Range("B41").Value = 'Project Documents'
$ws6.Range("C41").Value = 0.25
$ws6.Range("A48").Value = 'Dashboard query code'
$ws6.Range("B48").Value = 'Dashboard'
$ws6.Range("C48").Value = 0.5
$ws6.Range("A49").Value = 'Work on report appendix'
$ws6.Range("B49").Value = 'Report'
$ws6.Range("C49").Value = 1.5
$ws6.Range("A50").Value = 'Stand-up'
$ws6.Range("B50").Value = 'Internal meeting'
$ws6.Range("C50").Value = 0.25
$ws6.Range("A51").Value = 'Work on final ppt'
$ws6.Range("B51").Value = 'Final Presentation'
$ws6.Range("C51").Value = 2.75
$ws6.Range("A52").Value = 'Sprint planning meeting'
$ws6.Range("B52").Value = 'Client meeting'
$ws6.Range("C52").Value = 0.5
$ws6.Range("A53").Value = 'Post meeting communication with client'
$ws6.Range("B53").Value = 'Client communication'
$ws6.Range("C53").Value = 0.25
$ws6.Range("A54").Value = 'Discusssion with Ryan on Report Results'
$ws6.Range("B54").Value = 'Report'
$ws6.Range("C54").Value = 0.5
$ws6.Range("B55").Value = 'Code/Repo clean-up'
$ws6.Range("A55").Value = 'Work on cleaning up GitHub'
$ws6.Range("C55").Value = 2.5
$ws6.Range("A60").Value = 'Dashboard query code'
$ws6.Range("B60").Value = 'Dashboard'
$ws6.Range("C60").Value = 0.5
$ws6.Range("A61").Value = 'Manual data downloads for phase 2 testing'
$ws6.Range("B61").Value = 'Manual Data Downloads'
$ws6.Range("C61").Value = 2
$ws6.Range("A62").Value = 'Stand-up'
$ws6.Range("B62").Value = 'Internal Meetings'
$ws6.Range("C62").Value = 0.25
$ws6.Range("A63").Value = 'Discuss dashboard with Mitch'
$ws6.Range("B63").Value = 'Dashboard'
$ws6.Range("C63").Value = 0.25
$ws6.Range("A64").Value = 'Work on report and executive Summary'
$ws6.Range("B64").Value = 'Report'
$ws6.Range("C64").Value = 2.25
$ws6.Range("A71").Value = 'Working on dashboard'
$ws6.Range("B71").Value = 'Dashboard'
$ws6.Range("C71").Value = 3.5
$ws6.Range("A83").Value = 'Project organization'
$ws6.Range("B83").Value = 'Management'
$ws6.Range("C83").Value = 0.5
$ws6.Range("A82").Value = 'Working on report (results and formatting)'
$ws6.Range("B82").Value = 'Report'
$ws6.Range("C82").Value = 4.5

# --- Update view state: Week 5 loses tab-selection/active selection, Week 6 becomes active ---
$ws5.Activate()
$ws5.Range("O17").Select()

$ws6.Activate()
$ws6.Range("C61").Select()
